$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts existing B:E to C:F)
$ws.Columns.Item(2).Insert()

# Set header for new column
$ws.Cells.Item(1, 2).Value = "Variable"

# Set values for data rows
$ws.Cells.Item(2, 2).Value = "c1"
$ws.Cells.Item(3, 2).Value = "c1"
$ws.Cells.Item(4, 2).Value = "c1"
$ws.Cells.Item(5, 2).Value = "c1"

# Give the new column a best-fit-like width similar to the other data columns
$ws.Columns.Item(2).ColumnWidth = 6.8
